$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column F header ("Nível de Água (cm)") - reuse header style (same as
# A1:E1, i.e. bold/centered/wrapped with border) by copying format from E1.
# ---------------------------------------------------------------------------
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Nível de Água (cm)"

# ---------------------------------------------------------------------------
# Helper style cells already on the sheet that we reuse as format sources:
#   A2  -> date, border, center+vcenter+wrap  (numFmtId 14)
#   A10 -> border only, no alignment/number format (plain bordered cell)
# We paste-special (formats only) from these "donor" cells onto the new
# cells, then tweak alignment/number-format so the resulting style exactly
# matches the new xf records introduced by the edit (indices 6, 7, 8):
#   6 -> border + horizontal-center only (no vertical/wrap, General)
#   7 -> border + horizontal-center only, numFmtId 20 (time, "h:mm")
#   8 -> border + horizontal-center only, numFmtId 14 (date, "mm-dd-yy")
# ---------------------------------------------------------------------------

function Set-PlainCell($addr, $value) {
    $ws.Range("A10").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $text) {
    $ws.Range("A10").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).Value = "'" + $text
}

function Set-TimeCell($addr, $value) {
    $ws.Range("A10").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).NumberFormat = "h:mm"
    $ws.Range($addr).Value = $value
}

function Set-DateCell($addr, $value) {
    $ws.Range("A10").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).NumberFormat = "mm-dd-yy"
    $ws.Range($addr).Value = $value
}

function Set-DateWrapCell($addr, $value) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------------
# Row 2-3: extend existing rows with the new "Nível de Água (cm)" reading.
# ---------------------------------------------------------------------------
Set-PlainCell "F2" 75
Set-PlainCell "F3" 72

# ---------------------------------------------------------------------------
# Row 4: new measurement row (dates/times keep the original wrapped style
# used by rows 2-3 for column A, plain-centered for everything else).
# ---------------------------------------------------------------------------
Set-DateWrapCell "A4" 45803
Set-TimeCell      "B4" 0.59375
Set-TextCell       "C4" "25.7"
Set-PlainCell     "D4" 59
Set-PlainCell     "E4" 114
Set-PlainCell     "F4" 73

# ---------------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------------
Set-DateWrapCell "A5" 45803
Set-TimeCell      "B5" 0.59722222222222221
Set-TextCell       "C5" "24.6"
Set-PlainCell     "D5" 52
Set-PlainCell     "E5" 120
Set-PlainCell     "F5" 110

# ---------------------------------------------------------------------------
# Row 6: first row whose date cell uses the new plain-centered date style.
# ---------------------------------------------------------------------------
Set-DateCell      "A6" 45803
Set-TimeCell      "B6" 0.61111111111111116
Set-TextCell       "C6" "25.1"
Set-PlainCell     "D6" 48
Set-PlainCell     "E6" 115
Set-PlainCell     "F6" 135

# ---------------------------------------------------------------------------
# Row 7
# ---------------------------------------------------------------------------
Set-DateCell      "A7" 45803
Set-TimeCell      "B7" 0.62013888888888891
Set-TextCell       "C7" "25.6"
Set-PlainCell     "D7" 54
Set-PlainCell     "E7" 119
Set-PlainCell     "F7" 185

# ---------------------------------------------------------------------------
# Row 8
# ---------------------------------------------------------------------------
Set-DateCell      "A8" 45803
Set-TimeCell      "B8" 0.63402777777777775
Set-TextCell       "C8" "24.0"
Set-PlainCell     "D8" 53
Set-PlainCell     "E8" 115
Set-PlainCell     "F8" 164

# ---------------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------------
Set-DateCell      "A9" 45803
Set-TimeCell      "B9" 0.6430555555555556
Set-TextCell       "C9" "25.4"
Set-PlainCell     "D9" 57
Set-PlainCell     "E9" 118
Set-PlainCell     "F9" 127

# ---------------------------------------------------------------------------
# Row 10: extend the already-blank bordered row to column F.
# ---------------------------------------------------------------------------
$ws.Range("E10").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Selection moves to D14 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("D14").Select() | Out-Null
